# Updated cryptos list (price / 1h volume change refresh, plus a rank swap
# between Monero and BabyDogeCoin in rows 50-51), mirroring the scheduled
# GitHub Actions data refresh.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").Value = "68.481.52"
$ws.Range("E2").Value = "  +1.58%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "3.931.63"
$ws.Range("E3").Value = "  -0.08%  "

# Row 4 - TetherUSD (leading apostrophe forces text so "1.00" keeps its trailing zero)
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  -0.13%  "

# Row 5 - BNB
$ws.Range("D5").Value = "'487.73"
$ws.Range("E5").Value = "  +3.98%  "

# Row 6 - Solana
$ws.Range("D6").Value = "'148.74"
$ws.Range("E6").Value = "  +2.56%  "

# Row 7 - XRP (price unchanged)
$ws.Range("E7").Value = "  +0.95%  "

# Row 8 - USDC (price unchanged)
$ws.Range("E8").Value = "  -0.04%  "

# Row 9 - Cardano
$ws.Range("D9").Value = "'0.739"
$ws.Range("E9").Value = "  +0.72%  "

# Row 10 - Dogecoin (price unchanged)
$ws.Range("E10").Value = "  +1.32%  "

# Row 11 - ShibaInu
$ws.Range("D11").Value = "'0.0000356"
$ws.Range("E11").Value = "  +4.45%  "

# Row 12 - Avalanche
$ws.Range("D12").Value = "'43.15"
$ws.Range("E12").Value = "  -0.49%  "

# Row 13 - Polkadot
$ws.Range("D13").Value = "'10.75"
$ws.Range("E13").Value = "  +3.20%  "

# Row 14 - WrappedliquidstakedEther2.0
$ws.Range("D14").Value = "4.567.86"
$ws.Range("E14").Value = "  +0.55%  "

# Row 15 - Uniswap
$ws.Range("D15").Value = "'14.78"
$ws.Range("E15").Value = "  -1.52%  "

# Row 16 - WrappedEther
$ws.Range("D16").Value = "3.922.73"
$ws.Range("E16").Value = "  -0.08%  "

# Row 17 - TRON (price unchanged)
$ws.Range("E17").Value = "  -0.57%  "

# Row 18 - Chainlink
$ws.Range("D18").Value = "'20.06"
$ws.Range("E18").Value = "  +0.75%  "

# Row 19 - Polygon (price unchanged)
$ws.Range("E19").Value = "  -1.76%  "

# Row 20 - WrappedBTC
$ws.Range("D20").Value = "68.542.19"
$ws.Range("E20").Value = "  +1.39%  "

# Row 21 - BitcoinCash
$ws.Range("D21").Value = "'445.87"
$ws.Range("E21").Value = "  +2.75%  "

# Row 22 - ImmutableX (price unchanged)
$ws.Range("E22").Value = "  +4.29%  "

# Row 23 - InternetComputer(DFINITY)
$ws.Range("D23").Value = "'14.96"
$ws.Range("E23").Value = "  +1.76%  "

# Row 24 - Litecoin
$ws.Range("D24").Value = "'88.83"
$ws.Range("E24").Value = "  +1.01%  "

# Row 25 - RenderToken
$ws.Range("D25").Value = "'11.42"
$ws.Range("E25").Value = "  +19.28%  "

# Row 26 - Filecoin
$ws.Range("D26").Value = "'11.57"
$ws.Range("E26").Value = "  +12.66%  "

# Row 27 - PancakeSwap (price unchanged)
$ws.Range("E27").Value = "  +3.06%  "

# Row 28 - EthereumClassic
$ws.Range("D28").Value = "'39.00"
$ws.Range("E28").Value = "  +0.75%  "

# Row 29 - LEO (price unchanged)
$ws.Range("E29").Value = "  +1.51%  "

# Row 30 - Bittensor
$ws.Range("D30").Value = "'725.61"
$ws.Range("E30").Value = "  -0.08%  "

# Row 31 - Cosmos
$ws.Range("D31").Value = "'13.71"
$ws.Range("E31").Value = "  +0.80%  "

# Row 32 - Hedera: no change

# Row 33 - Toncoin
$ws.Range("D33").Value = "'2.90"
$ws.Range("E33").Value = "  +3.08%  "

# Row 34 - PEPE (price unchanged)
$ws.Range("E34").Value = "  +15.53%  "

# Row 35 - InjectiveProtocol
$ws.Range("D35").Value = "'42.48"
$ws.Range("E35").Value = "  -1.29%  "

# Row 36 - NEARProtocol
$ws.Range("D36").Value = "'6.20"
$ws.Range("E36").Value = "  +15.13%  "

# Row 37 - OKB
$ws.Range("D37").Value = "'61.12"
$ws.Range("E37").Value = "  +5.69%  "

# Row 38 - Kaspa (price unchanged)
$ws.Range("E38").Value = "  -3.37%  "

# Row 39 - TheGraph: no change

# Row 40 - Dai
$ws.Range("D40").Value = "'1.00"
$ws.Range("E40").Value = "  -0.01%  "

# Row 41 - Fetch.AI
$ws.Range("D41").Value = "'2.99"
$ws.Range("E41").Value = "  +14.85%  "

# Row 42 - VeChain
$ws.Range("D42").Value = "'0.0485"
$ws.Range("E42").Value = "  +1.27%  "

# Row 43 - ThetaToken
$ws.Range("D43").Value = "'3.19"
$ws.Range("E43").Value = "  +3.84%  "

# Row 44 - WEMIXToken: no change

# Row 45 - Stellar (price unchanged)
$ws.Range("E45").Value = "  +0.71%  "

# Row 46 - FirstDigitalUSD (price unchanged)
$ws.Range("E46").Value = "  -0.08%  "

# Row 47 - LidoDAOToken (price unchanged)
$ws.Range("E47").Value = "  +0.44%  "

# Row 48 - ARBITRUM (price unchanged)
$ws.Range("E48").Value = "  -0.99%  "

# Row 49 - ApeXProtocol (price unchanged)
$ws.Range("E49").Value = "  +1.24%  "

# Row 50 - was Monero, now BabyDogeCoin (rows 50/51 swapped places)
$ws.Range("B50").Value = "BabyDogeCoin"
$ws.Range("C50").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D50").Value = "0.0₆0345"
$ws.Range("E50").Value = "  +38.76%  "

# Row 51 - was BabyDogeCoin, now Monero
$ws.Range("B51").Value = "Monero"
$ws.Range("C51").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D51").Value = "'145.99"
$ws.Range("E51").Value = "  -0.05%  "
